$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "POR"
$ws.Range("C2").Value = 13.51428571428571
$ws.Range("B3").Value = "CLE"
$ws.Range("C3").Value = 13.01538461538462
$ws.Range("B4").Value = "DAL"
$ws.Range("C4").Value = 11.22105263157895
$ws.Range("B5").Value = "MIA"
$ws.Range("C5").Value = 12.18125
$ws.Range("B6").Value = "OKC"
$ws.Range("C6").Value = 12.33125
$ws.Range("B7").Value = "ATL"
$ws.Range("C7").Value = 10.98888888888889
$ws.Range("B8").Value = "WAS"
$ws.Range("C8").Value = 13.97857142857143
$ws.Range("B9").Value = "MIL"
$ws.Range("C9").Value = 12.08125
$ws.Range("B10").Value = "LAC"
$ws.Range("C10").Value = 13.07857142857143
$ws.Range("B11").Value = "SAS"
$ws.Range("C11").Value = 14.34705882352941
$ws.Range("B12").Value = "DET"
$ws.Range("C12").Value = 11.33571428571429
$ws.Range("B13").Value = "ORL"
$ws.Range("C13").Value = 12.61764705882353
$ws.Range("B14").Value = "UTA"
$ws.Range("C14").Value = 19.85
$ws.Range("B15").Value = "MEM"
$ws.Range("C15").Value = 13.8
$ws.Range("B16").Value = "HOU"
$ws.Range("C16").Value = 11.02222222222222
$ws.Range("B17").Value = "NOP"
$ws.Range("C17").Value = 12.44666666666667
$ws.Range("B18").Value = "DEN"
$ws.Range("C18").Value = 13.9875
$ws.Range("B19").Value = "LAL"
$ws.Range("C19").Value = 10.69411764705883
$ws.Range("B20").Value = "GSW"
$ws.Range("C20").Value = 13.30588235294118
$ws.Range("B21").Value = "IND"
$ws.Range("C21").Value = 13.98333333333333
$ws.Range("B22").Value = "CHO"
$ws.Range("C22").Value = 13.23125
$ws.Range("B23").Value = "CHI"
$ws.Range("C23").Value = 12.10666666666666
$ws.Range("B24").Value = "PHI"
$ws.Range("C24").Value = 11.575
$ws.Range("B25").Value = "BOS"
$ws.Range("C25").Value = 11.94444444444444
$ws.Range("B26").Value = "BRK"
$ws.Range("C26").Value = 12.21333333333333
$ws.Range("B27").Value = "TOR"
$ws.Range("C27").Value = 14.54
$ws.Range("B28").Value = "SAC"
$ws.Range("C28").Value = 14.17692307692308
$ws.Range("B29").Value = "PHO"
$ws.Range("C29").Value = 11.75
$ws.Range("B30").Value = "NYK"
$ws.Range("C30").Value = 10.7
$ws.Range("B31").Value = "MIN"
$ws.Range("C31").Value = 16.58571428571429
